# Fix dataset accession (associated_datasets, column CD) and year (column Y),
# plus trim a stray leading space from "last_author" (column AA).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("crispr_studies")

# --- Fix study_year typo on row 10 (20224 -> 2024) ---
$ws.Range("Y10").Value = 2024

# --- Trim leading space from last_author "Jonathan S. Weissman" ---
$ws.Range("AA2").Value = "Jonathan S. Weissman"
$ws.Range("AA3").Value = "Jonathan S. Weissman"
$ws.Range("AA11").Value = "Jonathan S. Weissman"
$ws.Range("AA12").Value = "Jonathan S. Weissman"

# --- Populate associated_datasets (column CD) with real dataset accessions ---
$ws.Range("CD2").Value = '[{"dataset_accession": "Gilbert LA (2014) - 1-PMID25307932", "dataset_uri": "https://orcs.thebiogrid.org/Download?type=screen&id=1161", "dataset_description": "Scores", "dataset_file_name": null}]'
$ws.Range("CD3").Value = '[{"dataset_accession": "Gilbert LA (2014) - 2-PMID25307932", "dataset_uri": "https://orcs.thebiogrid.org/Download?type=screen&id=5", "dataset_description": "Scores", "dataset_file_name": null}]'
$ws.Range("CD4").Value = '[{"dataset_accession": "Panea RI (2019) - 1-PMID31558468", "dataset_uri": "https://orcs.thebiogrid.org/Download?type=screen&id=2362", "dataset_description": "Scores", "dataset_file_name": null}]'
$ws.Range("CD5").Value = '[{"dataset_accession": "Panea RI (2019) - 2-PMID31558468", "dataset_uri": "https://orcs.thebiogrid.org/Download?type=screen&id=2363", "dataset_description": "Scores", "dataset_file_name": null}]'
$ws.Range("CD6").Value = '[{"dataset_accession": "Panea RI (2019) - 3-PMID31558468", "dataset_uri": "https://orcs.thebiogrid.org/Download?type=screen&id=2364", "dataset_description": "Scores", "dataset_file_name": null}]'
$ws.Range("CD7").Value = '[{"dataset_accession": "Panea RI (2019) - 4-PMID31558468", "dataset_uri": "https://orcs.thebiogrid.org/Download?type=screen&id=2365", "dataset_description": "Scores", "dataset_file_name": null}]'
$ws.Range("CD8").Value = '[{"dataset_accession": "Panea RI (2019) - 5-PMID31558468", "dataset_uri": "https://orcs.thebiogrid.org/Download?type=screen&id=2366", "dataset_description": "Scores", "dataset_file_name": null}]'
$ws.Range("CD9").Value = '[{"dataset_accession": "Panea RI (2019) - 6-PMID31558468", "dataset_uri": "https://orcs.thebiogrid.org/Download?type=screen&id=2367", "dataset_description": "Scores", "dataset_file_name": null}]'
$ws.Range("CD10").Value = '[{"dataset_accession": "PRJNA1170571", "dataset_uri": "https://www.ncbi.nlm.nih.gov/bioproject/PRJNA1170571", "dataset_description": "Raw sequence reads", "dataset_file_name": null},{"dataset_accession": "Chen Z (2024) - 1-PMID39547224", "dataset_uri": "https://orcs.thebiogrid.org/Download?type=screen&id=2373", "dataset_description": "Scores", "dataset_file_name": null}]'
$ws.Range("CD11").Value = '[{"dataset_accession": "Gilbert LA (2014) - 3-PMID25307932", "dataset_uri": "https://orcs.thebiogrid.org/Download?type=screen&id=6", "dataset_description": "Scores", "dataset_file_name": null}]'
$ws.Range("CD12").Value = '[{"dataset_accession": "Gilbert LA (2014) - 4-PMID25307932", "dataset_uri": "https://orcs.thebiogrid.org/Download?type=screen&id=1162", "dataset_description": "Scores", "dataset_file_name": null}]'

# --- Update the window view: freeze first column, scroll to Y1, select AA5 ---
$ws.Activate()
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("Y1").Select()
$ws.Range("AA5").Select()
